$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.5
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.9
$ws.Range("T2").Value = 1.33

# Row 3
$ws.Range("G3").Value = 6.25
$ws.Range("H3").Value = 3.8
$ws.Range("I3").Value = 1.55
$ws.Range("J3").Value = 6
$ws.Range("L3").Value = 2.1
$ws.Range("M3").Value = 1.05
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("S3").Value = 3.25
$ws.Range("W3").Value = 1.91
$ws.Range("X3").Value = 1.8
$ws.Range("AB3").Value = 67
$ws.Range("AD3").Value = 51
$ws.Range("AF3").Value = 7.5
$ws.Range("AI3").Value = 6.5
$ws.Range("AL3").Value = 11

# Row 4
$ws.Range("M4").Value = 1.04
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 3
$ws.Range("T4").Value = 1.36

# Row 5
$ws.Range("G5").Value = 1.06
$ws.Range("H5").Value = 7.6
$ws.Range("I5").Value = 29
$ws.Range("J5").Value = 1.28
$ws.Range("K5").Value = 3.35
$ws.Range("L5").Value = 19.5
$ws.Range("Q5").Value = 1.25
$ws.Range("R5").Value = 3.7
$ws.Range("S5").Value = 1.65
$ws.Range("T5").Value = 1.98
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 5.8
$ws.Range("AA5").Value = 11.5
$ws.Range("AB5").Value = 5.1
$ws.Range("AC5").Value = 10.25
$ws.Range("AD5").Value = 37
$ws.Range("AE5").Value = 20
$ws.Range("AF5").Value = 19
$ws.Range("AG5").Value = 45
$ws.Range("AH5").Value = 175
$ws.Range("AI5").Value = 100
$ws.Range("AK5").Value = 110
$ws.Range("AN5").Value = 300

# Row 7
$ws.Range("G7").Value = 1.1
$ws.Range("H7").Value = 9.5
$ws.Range("I7").Value = 23
$ws.Range("W7").Value = 2
$ws.Range("X7").Value = 1.75
$ws.Range("Y7").Value = 13
$ws.Range("AJ7").Value = 126
$ws.Range("AM7").Value = 151
$ws.Range("AN7").Value = 101
$ws.Range("AO7").Value = 351

# Row 8
$ws.Range("G8").Value = 2.45
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 2.63
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 2.3
$ws.Range("L8").Value = 3.1
$ws.Range("O8").Value = 1.18
$ws.Range("P8").Value = 4.5
$ws.Range("Q8").Value = 1.62
$ws.Range("R8").Value = 2.25
$ws.Range("S8").Value = 2.5
$ws.Range("T8").Value = 1.5
$ws.Range("U8").Value = 1.3
$ws.Range("V8").Value = 3.4
$ws.Range("W8").Value = 1.5
$ws.Range("X8").Value = 2.5
$ws.Range("Y8").Value = 12
$ws.Range("Z8").Value = 15
$ws.Range("AC8").Value = 17
$ws.Range("AE8").Value = 15
$ws.Range("AF8").Value = 7
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 34
$ws.Range("AI8").Value = 12
$ws.Range("AM8").Value = 19
$ws.Range("AN8").Value = 23
$ws.Range("AO8").Value = 101
$ws.Range("AP8").Value = 2.03
$ws.Range("AQ8").Value = 1.83

# Row 9
$ws.Range("G9").Value = 1.3
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 1.73
$ws.Range("K9").Value = 2.63
$ws.Range("O9").Value = 1.17
$ws.Range("P9").Value = 5
$ws.Range("Q9").Value = 1.53
$ws.Range("R9").Value = 2.4
$ws.Range("S9").Value = 2.25
$ws.Range("T9").Value = 1.57
$ws.Range("AA9").Value = 8.5
$ws.Range("AD9").Value = 23
$ws.Range("AF9").Value = 11
$ws.Range("AI9").Value = 23
$ws.Range("AK9").Value = 23
$ws.Range("AL9").Value = 101
$ws.Range("AP9").Value = 1.88
$ws.Range("AQ9").Value = 1.98
